$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = -3.6908300000000001
$ws.Range("K2").Value = 3.6051000000000002
$ws.Range("L2").Value = 0.72284099999999996

$ws.Range("J3").Value = 4.8814500000000001
$ws.Range("K3").Value = -1.55169
$ws.Range("L3").Value = -5.5954300000000003

$ws.Range("J4").Value = -7.3979299999999997
$ws.Range("K4").Value = 0.53811799999999999
$ws.Range("L4").Value = 7.1681400000000002

$ws.Range("J5").Value = 0.68991400000000003
$ws.Range("K5").Value = -4.21556
$ws.Range("L5").Value = 1.69537

$ws.Range("J6").Value = -2.8578199999999998
$ws.Range("K6").Value = -16.145800000000001
$ws.Range("L6").Value = 14.116300000000001

$ws.Range("J2:L7").Select()
